$p = $ppt.ActivePresentation

# --- 1) Handout master: datetimeFigureOut cached text 1/22/17 -> 1/23/17 ---
$hmDate = $p.HandoutMaster.HeadersFooters.DateAndTime
$hmDate.Text = "1/23/17"

# --- 2) Notes master: datetimeFigureOut cached text 1/22/17 -> 1/23/17 ---
$nmDate = $p.NotesMaster.HeadersFooters.DateAndTime
$nmDate.Text = "1/23/17"

# --- 3) Slide 9 speaker notes: merge the two runs of the second paragraph
#        ("5 " + "minutes. Please tell ...") into a single run of text,
#        leaving the first paragraph's wording untouched. Read the existing
#        paragraph text back (TextRange.Text already flattens a paragraph's
#        runs into one string) and rewrite the whole notes body with the
#        same paragraphs, so only the run split inside paragraph 2 changes. ---
$notesShape = $p.Slides.Item(9).NotesPage.Shapes.Item(2)
$notesTextRange = $notesShape.TextFrame.TextRange
$paraCount = $notesTextRange.Paragraphs().Count
$paraTexts = @()
for ($i = 1; $i -le $paraCount; $i++) {
    $paraTexts += $notesTextRange.Paragraphs($i, 1).Text
}
$notesTextRange.Text = [string]::Join([string][char]10, $paraTexts)
